$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversion del dia" note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 9.89 = 41642.53 pesos"), "1000 Bs = 9.91 = 41774.03 pesos"
$text = $text -replace [regex]::Escape("41642.53 pesos = 9.84 = 964.77 Bs"), "41774.03 pesos = 9.87 = 972.3 Bs"
$cellA1.Value2 = $text

# --- tasas: update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 100.9
$wsTasas.Range("O10").Value = 4215
$wsTasas.Range("N12").Value = 4233.27
$wsTasas.Range("O12").Value = 98.53
